# "Hores agost setembre" - fill in hours worked for the first half of
# September and clear out the placeholder dates that were never used.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C9 ("Treball a casa de nit") is no longer used once B9 is rewritten below,
# so clear it first -- this drops it from the shared-string table.
$ws.Range("C9").Value = $null

# 2016-09-07 / 2016-09-08: new entries ("4h" then "5h ").
$ws.Range("B10").Value = "4h"
$ws.Range("B11").Value = "5h "

# 2016-09-06: B9 used to read "3h centre 4h casa"; now it's simply "3h ".
$ws.Range("B9").Value = "3h "

# 2016-09-09 through 2016-09-15: remaining entries for the week.
$ws.Range("B12").Value = "2h"
$ws.Range("B13").Value = "?"
$ws.Range("B14").Value = "?"
$ws.Range("B15").Value = "?"
$ws.Range("B16").Value = "?"
$ws.Range("B17").Value = "?"
$ws.Range("B18").Value = "?"

# Rows 19-36 held placeholder dates through mid-October that were never
# filled in; clear the date values but keep the row formatting in place.
for ($row = 19; $row -le 36; $row++) {
    $ws.Cells.Item($row, 1).Value = $null
}

# Restore the active cell/selection recorded in the workbook view.
$ws.Range("B14").Select()
